# Team-Meeting-2.docx regeneration edit
#
# 1) Drop the three leading "site nav" paragraphs (Home / Back to Home /
#    Download Word Document hyperlinks) that used to sit before the
#    document's opening bookmark.
# 2) Give every table in the document an explicit 100% preferred width
#    (percentage-based) instead of the old "auto" width, matching the
#    regenerated table formatting used across the site.

$d = $word.ActiveDocument

# --- Step 1: remove the first three paragraphs -----------------------
$firstPara = $d.Paragraphs.Item(1)
$thirdPara = $d.Paragraphs.Item(3)
$navRange = $d.Range($firstPara.Range.Start, $thirdPara.Range.End)
$navRange.Delete()

# --- Step 2: set every table to a 100% preferred width ----------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $tbl.PreferredWidthType = 2   # wdPreferredWidthPercent
    $tbl.PreferredWidth = 250     # -> serializes as <w:tblW w:type="pct" w:w="5000"/> (100%)
}
